# Fill in the empty "Заняття"/"Здача" date cells for Тема 5
# (ЛР12, ЛР13, Л14, ЛР14, ЛР15 rows) of the schedule table with the
# dates added by the commit.
#
# Table layout (1-based, as used by Table.Cell(row, col)):
#   col 1 = №№ (numbering)
#   col 2 = ІПЗ-21 Заняття     col 3 = ІПЗ-21 Здача
#   col 4 = ІПЗ-22 Заняття     col 5 = ІПЗ-22 Здача
#   col 6 = Тема (code)        col 7 = Тема (description)

$d = $word.ActiveDocument
$t = $d.Tables(1)

function Set-CellDate($table, $row, $col, $text) {
    $rng = $table.Cell($row, $col).Range
    $rng.Text = $text
    $rng.Font.Name = "Times New Roman"
    $rng.Font.Size = 14
    $rng.Font.SizeBi = 14
}

# Row 36 (ЛР12): ІПЗ-22 Заняття/Здача were empty
Set-CellDate $t 36 4 "04.05"
Set-CellDate $t 36 5 "07.05"

# Row 37 (ЛР13): all four date cells were empty
Set-CellDate $t 37 2 "06.05"
Set-CellDate $t 37 3 "10.05"
Set-CellDate $t 37 4 "05.05"
Set-CellDate $t 37 5 "10.05"

# Row 38 (Л14): ІПЗ-21 Заняття and ІПЗ-22 Заняття were empty
Set-CellDate $t 38 2 "07.05"
Set-CellDate $t 38 4 "05.05"

# Row 39 (ЛР14): all four date cells were empty
Set-CellDate $t 39 2 "07.05"
Set-CellDate $t 39 3 "11.05"
Set-CellDate $t 39 4 "07.05"
Set-CellDate $t 39 5 "11.05"

# Row 40 (ЛР15): ІПЗ-22 Заняття/Здача were empty
Set-CellDate $t 40 4 "07.05"
Set-CellDate $t 40 5 "12.05"

Write-Output "done"
